$wb = $excel.ActiveWorkbook

# Remove the "TODO" worksheet entirely - all of its open TODOs have been
# handled / folded into the "Fragen an Ivan" sheet.
$todoSheet = $wb.Worksheets.Item("TODO")
[void]$todoSheet.Delete()

# Append the remaining newly-raised questions to the "Fragen an Ivan" sheet.
$ws = $wb.Worksheets.Item("Fragen an Ivan")

$ws.Range("A18").Value = "Richtig, dass Zeichnung und Plan StillImageRepresentation? Auch wenn PDF?"
$ws.Range("A19").Value = "Lage alles 0-1 und denne mehreri Lage Objekt, Lage nur für d Abkürzig vom ganze, wie verlinke? -> Ivan "
$ws.Range("A20").Value = "Sollten min und max überhaupt definiert werden? "

# Make this the active sheet/cell, matching the post-edit selection state.
$ws.Activate()
[void]$ws.Range("A20").Select()
